# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" right before the "总计" (Total) sheet,
#   populated with per-fund holding detail (same shape as the other
#   quarterly sheets).
# - Prepend a "2022-Q1" row to the "总计" summary sheet (date / count /
#   market value), shifting the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Create the new "2022-Q1" sheet immediately before "总计" and rename it.
#
# NOTE: this COM layer hands out position-based handles — any handle
# obtained before a sheet is inserted/removed can silently start
# pointing at a different sheet once the sheet list shifts. So we only
# use $wsTotal / $wsQ1 to perform the Add() + Name assignment, then
# throw them away and re-resolve every worksheet *by name* afterwards.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Add($wsTotal)
$wsQ1.Name = "2022-Q1"

# ---- fresh, stable references for everything that follows ----
$wsStyleDonor = $wb.Worksheets.Item("2021-Q1")   # any detail sheet; carries the shared "s=2" bold/border style
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("总计")

# ---- copy the shared header / index style (cellXfs index 2) over ----
$wsStyleDonor.Range("B1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$wsStyleDonor.Range("A2").Copy()
$wsQ1.Range("A2:A11").PasteSpecial(-4122)  # xlPasteFormats

# ---- header row ----
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# ---- fund detail rows ----
# columns: A=index  B=code(text)  C=name(text)  D=scale(text)
#          E=stock position(text)  F=position share(text)
#          G=market value(text)  H=rank(number)
$q1Rows = @(
    @(0, "011866", "广发价值增长混合型证券投资基金A",     "18.06", "92.01", "4.78", "0.8633", 9),
    @(1, "002624", "广发优企精选灵活配置混合A",           "12.98", "92.40", "5.68", "0.7373", 9),
    @(2, "270025", "广发行业领先混合A",                   "11.11", "91.67", "4.60", "0.5111", 10),
    @(3, "960001", "广发行业领先混合H",                   "11.11", "91.67", "4.60", "0.5111", 10),
    @(4, "501070", "广发睿阳三年定期开放混合",             "7.06",  "50.14", "5.08", "0.3586", 3),
    @(5, "210002", "金鹰红利价值混合",                     "0.88",  "69.18", "5.44", "0.0479", 2),
    @(6, "011867", "广发价值增长混合型证券投资基金C",     "0.81",  "92.01", "4.78", "0.0387", 9),
    @(7, "008353", "泰达宏利消费行业量化精选混合A",       "0.52",  "92.25", "3.15", "0.0164", 4),
    @(8, "010021", "广发优企精选灵活配置混合C",           "0.15",  "92.40", "5.68", "0.0085", 9),
    @(9, "008354", "泰达宏利消费行业量化精选混合C",       "0.12",  "92.25", "3.15", "0.0038", 4)
)

for ($i = 0; $i -lt $q1Rows.Length; $i++) {
    $r = $i + 2
    $row = $q1Rows[$i]

    $wsQ1.Cells.Item($r, 1).Value = $row[0]

    # text-like columns: force Text format first so numeric-looking
    # strings (fund codes, "18.06", …) are NOT coerced into numbers,
    # then reset the style back to Normal so no stray "s" attribute
    # (and no extra numFmt) leaks into the saved XML.
    for ($c = 2; $c -le 7; $c++) {
        $cell = $wsQ1.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
        $cell.Style = "Normal"
    }

    $wsQ1.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# Prepend the 2022-Q1 summary row to "总计", pushing the rest down.
# ---------------------------------------------------------------------
$totalRows = @(
    @("2022-Q1", 10, 3.1),
    @("2021-Q4", 9, 4.9),
    @("2021-Q3", 5, 1.14),
    @("2021-Q2", 4, 2.13),
    @("2021-Q1", 3, 0.78),
    @("2020-Q4", 1, 0.66)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $wsTotal.Cells.Item($r, 1).Value = $i
    $wsTotal.Cells.Item($r, 2).Value = $row[0]
    $wsTotal.Cells.Item($r, 3).Value = $row[1]
    $wsTotal.Cells.Item($r, 4).Value = $row[2]
}
